# Updated cryptos list refresh.
# Applies the new Price / Volume(1h) figures scraped for this run, and
# fixes the Aptos/USDe and Cosmos/dogwifhat row ordering (ranking swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $text) {
    # Force the literal text into the cell (Excel otherwise auto-coerces
    # plain-looking numerics like "148.00" or "0.484" into real numbers,
    # which would silently drop the trailing zeros / exact formatting the
    # source site renders). Stamping NumberFormat "@" first makes the COM
    # layer store it as text; resetting the style back to Normal afterwards
    # keeps the cell's style identical to the untouched cells around it.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Row 38/39 swap: Aptos now ranks above USDe ---------------------------
Set-TextCell $ws "B38" "Aptos"
Set-TextCell $ws "C38" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws "D38" "7.93"
Set-TextCell $ws "E38" "  -6.58%  "

Set-TextCell $ws "B39" "USDe"
Set-TextCell $ws "C39" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws "D39" "1.00"
Set-TextCell $ws "E39" "  +0.00%  "

# --- Row 49/50 swap: dogwifhat now ranks above Cosmos ----------------------
Set-TextCell $ws "B49" "dogwifhat"
Set-TextCell $ws "C49" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell $ws "D49" "2.47"
Set-TextCell $ws "E49" "  -10.97%  "

Set-TextCell $ws "B50" "Cosmos"
Set-TextCell $ws "C50" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D50" "7.48"
Set-TextCell $ws "E50" "  -4.56%  "

# --- Price / Volume(1h) refresh for every other row -------------------------
$updates = @{
    2  = @{ D = "66.789.25";  E = "  -4.16%  " }
    3  = @{ D = "3.462.23";   E = "  -4.14%  " }
    4  = @{ E = "  -0.02%  " }
    5  = @{ D = "604.02";     E = "  -4.28%  " }
    6  = @{ D = "148.00";     E = "  -6.68%  " }
    7  = @{ D = "3.458.60";   E = "  -4.21%  " }
    8  = @{ E = "  +0.04%  " }
    9  = @{ D = "0.484";      E = "  -2.25%  " }
    10 = @{ E = "  -4.82%  " }
    11 = @{ D = "7.50";       E = "  -0.95%  " }
    12 = @{ E = "  -4.02%  " }
    13 = @{ D = "0.0000215";  E = "  -5.34%  " }
    14 = @{ D = "31.76";      E = "  -5.86%  " }
    15 = @{ D = "4.044.02";   E = "  -4.24%  " }
    16 = @{ D = "3.454.87";   E = "  -4.50%  " }
    17 = @{ D = "66.839.07";  E = "  -3.95%  " }
    19 = @{ D = "6.45";       E = "  -4.44%  " }
    20 = @{ D = "15.38";      E = "  -4.79%  " }
    21 = @{ D = "10.01";      E = "  -2.76%  " }
    22 = @{ D = "440.74";     E = "  -4.89%  " }
    23 = @{ D = "0.610";      E = "  -5.84%  " }
    24 = @{ D = "78.65";      E = "  -0.43%  " }
    25 = @{ E = "  -0.08%  " }
    26 = @{ D = "3.597.65";   E = "  -4.23%  " }
    27 = @{ E = "  -10.33%  " }
    28 = @{ D = "9.88";       E = "  -8.24%  " }
    29 = @{ D = "8.44";       E = "  -10.20%  " }
    30 = @{ E = "  -6.71%  " }
    31 = @{ D = "1.61";       E = "  -6.87%  " }
    32 = @{ D = "0.168";      E = "  -2.86%  " }
    33 = @{ E = "  +0.03%  " }
    34 = @{ D = "25.46";      E = "  -4.23%  " }
    35 = @{ D = "6.09";       E = "  -7.46%  " }
    36 = @{ D = "3.451.77";   E = "  -4.37%  " }
    37 = @{ E = "  -7.75%  " }
    40 = @{ D = "0.998";      E = "  -0.14%  " }
    41 = @{ D = "173.74";     E = "  -2.97%  " }
    42 = @{ D = "0.0892";     E = "  -3.75%  " }
    43 = @{ D = "2.17";       E = "  -11.13%  " }
    44 = @{ D = "5.42";       E = "  -5.00%  " }
    45 = @{ D = "0.884";      E = "  -3.31%  " }
    46 = @{ D = "29.11";      E = "  -9.04%  " }
    47 = @{ D = "45.95";      E = "  -0.14%  " }
    48 = @{ D = "1.25";       E = "  -10.56%  " }
    51 = @{ D = "0.986";      E = "  -5.86%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        Set-TextCell $ws "D$row" $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        Set-TextCell $ws "E$row" $vals["E"]
    }
}
